$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743717012329/Ggg0CQFNgb-camera-video-segment-driver-1743717014829.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIRRJK7UM%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T180033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDoFSR9MeS564y4q87lBejrW%2B7pd%2Fccj8z%2FWn4l8KMvXwIgL7%2FgjwG%2Bo%2BULh5cqhhLaDwYiIPhEU1AwnloKFWRcu2Eq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDDSOLoFcGjDwjBgVDyq6A%2FzJ%2Bozcxm5L%2FpVLgvQg0LvX6DpViGfE4dF%2FQmLlluJBYx2fhZCHRpvlc0dH5%2F59%2FhCEMyADm%2Fr8CYZECA7l%2B352DCWKduT95cksa8HMWYsUpfyjljUV1GmurcmPr3Xyvt4Qvdx%2Fdb3C1ltHNSB3GhuFAMZf7E%2FueYEb1OdaWpI%2Bhf1fSXeyszMZgB5zI3bhxXURXvNf4PAc1C0wIn8YDd7uOkGgSt31FjFFG15zye4qvcxqExblGqPHqlsm46hg5M2x4ftvFDw9xLhmDqfeNzqtzun7rfOvkkVxZlNVtixzIb%2Bb8ts8dCaB68UR93x2GDq6C%2FviYquWZ%2BOhFXTe%2FQ9R7wYCgUSvPjuYAd7WcWUi3bWsvKkirucozE5QV8aql6KsHEhZDNEh%2BV3fBjOMhlzO7JwOZwYIiOZ0ts1Lytwrabz7NZTImw9RTlObvrJK9zYKutcOkMq63lnSyG4lI5ktbbxDkblUlVAGOuD44Dvzq8RiLYLMcmKvYkobLH3i2V%2FnjIWA7JI%2F%2FDGF51iGONbtUB8XVHfol6cf5CrdT7NSTfBIAulekBO8iNCCBDdZYcU4HyUmhvWbrNQw3aK%2FvwY6pQEzgTXwr2F2dmiLBe2tam1xrRi04gmviYTfozK5oGvDm74Nx1gCrM3qub49XdX%2B%2FUCyxNHQcecE%2BMNWoVpAGOSbzPalbaEBM6H84WN%2FMYM8Wor0fUgOCR4pE8qppzl6ZebWkrYww%2BP9DemhW8dAGd0evQyUBNt8o7HaCctSRM%2BbNuF3n0GFnFS0ily1MFM7GjsJpNOZLmpJMQF2ZOmal2xaINNtXys%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2002%3A00%3A33%20GMT&X-Amz-Signature=91783fcf75c27b358987766eb4240e8602efda6def0d391165fc46e963e37b7b'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743713574485/FgwSeOefdQ-camera-video-segment-driver-1743713576985.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIRRJK7UM%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T180033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDoFSR9MeS564y4q87lBejrW%2B7pd%2Fccj8z%2FWn4l8KMvXwIgL7%2FgjwG%2Bo%2BULh5cqhhLaDwYiIPhEU1AwnloKFWRcu2Eq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDDSOLoFcGjDwjBgVDyq6A%2FzJ%2Bozcxm5L%2FpVLgvQg0LvX6DpViGfE4dF%2FQmLlluJBYx2fhZCHRpvlc0dH5%2F59%2FhCEMyADm%2Fr8CYZECA7l%2B352DCWKduT95cksa8HMWYsUpfyjljUV1GmurcmPr3Xyvt4Qvdx%2Fdb3C1ltHNSB3GhuFAMZf7E%2FueYEb1OdaWpI%2Bhf1fSXeyszMZgB5zI3bhxXURXvNf4PAc1C0wIn8YDd7uOkGgSt31FjFFG15zye4qvcxqExblGqPHqlsm46hg5M2x4ftvFDw9xLhmDqfeNzqtzun7rfOvkkVxZlNVtixzIb%2Bb8ts8dCaB68UR93x2GDq6C%2FviYquWZ%2BOhFXTe%2FQ9R7wYCgUSvPjuYAd7WcWUi3bWsvKkirucozE5QV8aql6KsHEhZDNEh%2BV3fBjOMhlzO7JwOZwYIiOZ0ts1Lytwrabz7NZTImw9RTlObvrJK9zYKutcOkMq63lnSyG4lI5ktbbxDkblUlVAGOuD44Dvzq8RiLYLMcmKvYkobLH3i2V%2FnjIWA7JI%2F%2FDGF51iGONbtUB8XVHfol6cf5CrdT7NSTfBIAulekBO8iNCCBDdZYcU4HyUmhvWbrNQw3aK%2FvwY6pQEzgTXwr2F2dmiLBe2tam1xrRi04gmviYTfozK5oGvDm74Nx1gCrM3qub49XdX%2B%2FUCyxNHQcecE%2BMNWoVpAGOSbzPalbaEBM6H84WN%2FMYM8Wor0fUgOCR4pE8qppzl6ZebWkrYww%2BP9DemhW8dAGd0evQyUBNt8o7HaCctSRM%2BbNuF3n0GFnFS0ily1MFM7GjsJpNOZLmpJMQF2ZOmal2xaINNtXys%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2002%3A00%3A33%20GMT&X-Amz-Signature=d287ec2551dc3fa93ded34f719744b8c0b48e32de6484655eaffc551163d6ad3'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712999664/UzzO0YQTJX-camera-video-segment-driver-1743713002164.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIRRJK7UM%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T180033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDoFSR9MeS564y4q87lBejrW%2B7pd%2Fccj8z%2FWn4l8KMvXwIgL7%2FgjwG%2Bo%2BULh5cqhhLaDwYiIPhEU1AwnloKFWRcu2Eq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDDSOLoFcGjDwjBgVDyq6A%2FzJ%2Bozcxm5L%2FpVLgvQg0LvX6DpViGfE4dF%2FQmLlluJBYx2fhZCHRpvlc0dH5%2F59%2FhCEMyADm%2Fr8CYZECA7l%2B352DCWKduT95cksa8HMWYsUpfyjljUV1GmurcmPr3Xyvt4Qvdx%2Fdb3C1ltHNSB3GhuFAMZf7E%2FueYEb1OdaWpI%2Bhf1fSXeyszMZgB5zI3bhxXURXvNf4PAc1C0wIn8YDd7uOkGgSt31FjFFG15zye4qvcxqExblGqPHqlsm46hg5M2x4ftvFDw9xLhmDqfeNzqtzun7rfOvkkVxZlNVtixzIb%2Bb8ts8dCaB68UR93x2GDq6C%2FviYquWZ%2BOhFXTe%2FQ9R7wYCgUSvPjuYAd7WcWUi3bWsvKkirucozE5QV8aql6KsHEhZDNEh%2BV3fBjOMhlzO7JwOZwYIiOZ0ts1Lytwrabz7NZTImw9RTlObvrJK9zYKutcOkMq63lnSyG4lI5ktbbxDkblUlVAGOuD44Dvzq8RiLYLMcmKvYkobLH3i2V%2FnjIWA7JI%2F%2FDGF51iGONbtUB8XVHfol6cf5CrdT7NSTfBIAulekBO8iNCCBDdZYcU4HyUmhvWbrNQw3aK%2FvwY6pQEzgTXwr2F2dmiLBe2tam1xrRi04gmviYTfozK5oGvDm74Nx1gCrM3qub49XdX%2B%2FUCyxNHQcecE%2BMNWoVpAGOSbzPalbaEBM6H84WN%2FMYM8Wor0fUgOCR4pE8qppzl6ZebWkrYww%2BP9DemhW8dAGd0evQyUBNt8o7HaCctSRM%2BbNuF3n0GFnFS0ily1MFM7GjsJpNOZLmpJMQF2ZOmal2xaINNtXys%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2002%3A00%3A33%20GMT&X-Amz-Signature=218f6890cb48c50349f68137f93b192775fa96bf8c9e67a666dd42d85ff5c545'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712546260/Jh2p7rn04n-camera-video-segment-driver-1743712548760.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIRRJK7UM%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T180033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDoFSR9MeS564y4q87lBejrW%2B7pd%2Fccj8z%2FWn4l8KMvXwIgL7%2FgjwG%2Bo%2BULh5cqhhLaDwYiIPhEU1AwnloKFWRcu2Eq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDDSOLoFcGjDwjBgVDyq6A%2FzJ%2Bozcxm5L%2FpVLgvQg0LvX6DpViGfE4dF%2FQmLlluJBYx2fhZCHRpvlc0dH5%2F59%2FhCEMyADm%2Fr8CYZECA7l%2B352DCWKduT95cksa8HMWYsUpfyjljUV1GmurcmPr3Xyvt4Qvdx%2Fdb3C1ltHNSB3GhuFAMZf7E%2FueYEb1OdaWpI%2Bhf1fSXeyszMZgB5zI3bhxXURXvNf4PAc1C0wIn8YDd7uOkGgSt31FjFFG15zye4qvcxqExblGqPHqlsm46hg5M2x4ftvFDw9xLhmDqfeNzqtzun7rfOvkkVxZlNVtixzIb%2Bb8ts8dCaB68UR93x2GDq6C%2FviYquWZ%2BOhFXTe%2FQ9R7wYCgUSvPjuYAd7WcWUi3bWsvKkirucozE5QV8aql6KsHEhZDNEh%2BV3fBjOMhlzO7JwOZwYIiOZ0ts1Lytwrabz7NZTImw9RTlObvrJK9zYKutcOkMq63lnSyG4lI5ktbbxDkblUlVAGOuD44Dvzq8RiLYLMcmKvYkobLH3i2V%2FnjIWA7JI%2F%2FDGF51iGONbtUB8XVHfol6cf5CrdT7NSTfBIAulekBO8iNCCBDdZYcU4HyUmhvWbrNQw3aK%2FvwY6pQEzgTXwr2F2dmiLBe2tam1xrRi04gmviYTfozK5oGvDm74Nx1gCrM3qub49XdX%2B%2FUCyxNHQcecE%2BMNWoVpAGOSbzPalbaEBM6H84WN%2FMYM8Wor0fUgOCR4pE8qppzl6ZebWkrYww%2BP9DemhW8dAGd0evQyUBNt8o7HaCctSRM%2BbNuF3n0GFnFS0ily1MFM7GjsJpNOZLmpJMQF2ZOmal2xaINNtXys%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2002%3A00%3A33%20GMT&X-Amz-Signature=769eea399cdeb43aa05d661859d38afcece60ba0d7cd3c45194a4d000898950c'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743703000742/spTbAbfADq-camera-video-segment-driver-1743703003242.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIRRJK7UM%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T180033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJ3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDoFSR9MeS564y4q87lBejrW%2B7pd%2Fccj8z%2FWn4l8KMvXwIgL7%2FgjwG%2Bo%2BULh5cqhhLaDwYiIPhEU1AwnloKFWRcu2Eq3QMIFhAEGgw3ODEyMDQ5NDIyNDQiDDSOLoFcGjDwjBgVDyq6A%2FzJ%2Bozcxm5L%2FpVLgvQg0LvX6DpViGfE4dF%2FQmLlluJBYx2fhZCHRpvlc0dH5%2F59%2FhCEMyADm%2Fr8CYZECA7l%2B352DCWKduT95cksa8HMWYsUpfyjljUV1GmurcmPr3Xyvt4Qvdx%2Fdb3C1ltHNSB3GhuFAMZf7E%2FueYEb1OdaWpI%2Bhf1fSXeyszMZgB5zI3bhxXURXvNf4PAc1C0wIn8YDd7uOkGgSt31FjFFG15zye4qvcxqExblGqPHqlsm46hg5M2x4ftvFDw9xLhmDqfeNzqtzun7rfOvkkVxZlNVtixzIb%2Bb8ts8dCaB68UR93x2GDq6C%2FviYquWZ%2BOhFXTe%2FQ9R7wYCgUSvPjuYAd7WcWUi3bWsvKkirucozE5QV8aql6KsHEhZDNEh%2BV3fBjOMhlzO7JwOZwYIiOZ0ts1Lytwrabz7NZTImw9RTlObvrJK9zYKutcOkMq63lnSyG4lI5ktbbxDkblUlVAGOuD44Dvzq8RiLYLMcmKvYkobLH3i2V%2FnjIWA7JI%2F%2FDGF51iGONbtUB8XVHfol6cf5CrdT7NSTfBIAulekBO8iNCCBDdZYcU4HyUmhvWbrNQw3aK%2FvwY6pQEzgTXwr2F2dmiLBe2tam1xrRi04gmviYTfozK5oGvDm74Nx1gCrM3qub49XdX%2B%2FUCyxNHQcecE%2BMNWoVpAGOSbzPalbaEBM6H84WN%2FMYM8Wor0fUgOCR4pE8qppzl6ZebWkrYww%2BP9DemhW8dAGd0evQyUBNt8o7HaCctSRM%2BbNuF3n0GFnFS0ily1MFM7GjsJpNOZLmpJMQF2ZOmal2xaINNtXys%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2002%3A00%3A33%20GMT&X-Amz-Signature=fe82ae92290d2d5fe2908454d4c1422cb707d42cae6035696e528219144d40c5'
